# Updates cryptos list (price/volume columns D & E, and a couple of
# coin-name/link swaps in B & C) to match the Thu Jul  4 22:27:18 UTC 2024
# GitHub Actions refresh of cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to plain text first so Excel's COM layer
# doesn't auto-coerce numeric-looking strings (e.g. "0.446", "1.00") into
# real numbers -- the source data must stay as inline/shared text, exactly
# as it was before this edit.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '58.435.61'
$ws.Range("E2").Value = '  -2.92%  '
# Row 3
$ws.Range("D3").Value = '3.150.05'
$ws.Range("E3").Value = '  -4.41%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
$ws.Range("D5").Value = '526.87'
$ws.Range("E5").Value = '  -5.33%  '
# Row 6
$ws.Range("D6").Value = '133.81'
$ws.Range("E6").Value = '  -5.03%  '
# Row 7
$ws.Range("E7").Value = '  -0.08%  '
# Row 8
$ws.Range("D8").Value = '3.150.03'
$ws.Range("E8").Value = '  -4.45%  '
# Row 9
$ws.Range("D9").Value = '0.446'
$ws.Range("E9").Value = '  -4.47%  '
# Row 10
$ws.Range("D10").Value = '7.32'
$ws.Range("E10").Value = '  -7.25%  '
# Row 11
$ws.Range("D11").Value = '0.109'
$ws.Range("E11").Value = '  -7.65%  '
# Row 12
$ws.Range("D12").Value = '0.375'
$ws.Range("E12").Value = '  -8.03%  '
# Row 13
$ws.Range("D13").Value = '3.695.08'
$ws.Range("E13").Value = '  -4.33%  '
# Row 14
$ws.Range("E14").Value = '  -0.33%  '
# Row 15
$ws.Range("D15").Value = '25.29'
$ws.Range("E15").Value = '  -5.05%  '
# Row 16
$ws.Range("D16").Value = '3.160.56'
$ws.Range("E16").Value = '  -4.23%  '
# Row 17
$ws.Range("D17").Value = '58.411.92'
$ws.Range("E17").Value = '  -2.99%  '
# Row 18
$ws.Range("D18").Value = '0.0000152'
$ws.Range("E18").Value = '  -6.88%  '
# Row 19
$ws.Range("D19").Value = '5.76'
$ws.Range("E19").Value = '  -4.66%  '
# Row 20
$ws.Range("D20").Value = '13.04'
$ws.Range("E20").Value = '  -4.68%  '
# Row 21
$ws.Range("D21").Value = '7.92'
$ws.Range("E21").Value = '  -6.80%  '
# Row 22
$ws.Range("D22").Value = '343.74'
$ws.Range("E22").Value = '  -7.82%  '
# Row 23
$ws.Range("E23").Value = '  +0.01%  '
# Row 24
$ws.Range("D24").Value = '0.511'
$ws.Range("E24").Value = '  -3.74%  '
# Row 25
$ws.Range("D25").Value = '67.48'
$ws.Range("E25").Value = '  -7.21%  '
# Row 26
$ws.Range("D26").Value = '3.293.17'
$ws.Range("E26").Value = '  -3.99%  '
# Row 27
$ws.Range("D27").Value = '0.172'
$ws.Range("E27").Value = '  -0.91%  '
# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.20%  '
# Row 29
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0' + ([string][char]0x2083) + '0940'
$ws.Range("E29").Value = '  -7.78%  '
# Row 30
$ws.Range("D30").Value = '6.81'
$ws.Range("E30").Value = '  -2.83%  '
# Row 31
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.07%  '
# Row 32
$ws.Range("E32").Value = '  +4.78%  '
# Row 33
$ws.Range("D33").Value = '1.87'
$ws.Range("E33").Value = '  -7.16%  '
# Row 34
$ws.Range("D34").Value = '6.89'
$ws.Range("E34").Value = '  -7.14%  '
# Row 35
$ws.Range("D35").Value = '21.49'
$ws.Range("E35").Value = '  -4.56%  '
# Row 36
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '159.82'
$ws.Range("E36").Value = '  -3.63%  '
# Row 37
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '4.83'
$ws.Range("E37").Value = '  -3.78%  '
# Row 38
$ws.Range("D38").Value = '6.24'
$ws.Range("E38").Value = '  -5.54%  '
# Row 39
$ws.Range("D39").Value = '1.38'
$ws.Range("E39").Value = '  -8.83%  '
# Row 40
$ws.Range("D40").Value = '0.0687'
$ws.Range("E40").Value = '  -4.46%  '
# Row 41
$ws.Range("D41").Value = '3.184.71'
$ws.Range("E41").Value = '  -4.35%  '
# Row 42
$ws.Range("D42").Value = '40.43'
$ws.Range("E42").Value = '  -2.84%  '
# Row 43
$ws.Range("D43").Value = '23.88'
$ws.Range("E43").Value = '  -6.12%  '
# Row 44
$ws.Range("D44").Value = '0.695'
$ws.Range("E44").Value = '  -7.03%  '
# Row 45
$ws.Range("D45").Value = '1.09'
$ws.Range("E45").Value = '  -2.16%  '
# Row 46
$ws.Range("D46").Value = '3.94'
$ws.Range("E46").Value = '  -3.36%  '
# Row 47
$ws.Range("E47").Value = '  +0.02%  '
# Row 48
$ws.Range("D48").Value = '1.46'
$ws.Range("E48").Value = '  -6.35%  '
# Row 49
$ws.Range("D49").Value = '2.292.39'
$ws.Range("E49").Value = '  -1.04%  '
# Row 50
$ws.Range("D50").Value = '6.18'
$ws.Range("E50").Value = '  -2.19%  '
# Row 51
$ws.Range("D51").Value = '20.57'
$ws.Range("E51").Value = '  -4.01%  '

# Restore the default (Normal) style on those cells so no residual
# number-format override is left behind on the saved worksheet.
$ws.Range("D2:E51").Style = "Normal"
